# Wells_Missing_Data.xlsx update
# "Added meter offsets and welltop to ground measurements
#  based on field visit to Sagehen on 10/18/2020"
#
# Column A lists wells still missing "top2gnd" measurements. After the
# 10/18/2020 field visit, three wells (KHR-1, EFF-XA1N, EET-2) now have
# their top-of-ground measurement, so they are removed from the "still
# missing" list in column A: the remaining entries below them shift up,
# and the now-unused rows at the bottom of the list are cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the three "resolved" well IDs up into A2:A4 (they used to live
# at A5, A10 and A11 respectively).
$ws.Range("A2").Value = "KHR-1"
$ws.Range("A3").Value = "EFF-XA1N"
$ws.Range("A4").Value = "EET-2"

# These rows keep their formatting/border but no longer hold a value.
$ws.Range("A6").ClearContents()
$ws.Range("A7").ClearContents()
$ws.Range("A8").ClearContents()
$ws.Range("A13").ClearContents()
$ws.Range("A14").ClearContents()
$ws.Range("A15").ClearContents()

# These rows are fully emptied out (no leftover cell/format at all),
# matching the now-shorter "missing data" list.
$ws.Range("A5").Clear()
$ws.Range("A10").Clear()
$ws.Range("A11").Clear()

# Update the current selection to reflect where the editor ended up.
$ws.Range("C11").Select()
